$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 134
$ws.Range("F5").Value = 368
$ws.Range("F6").Value = 765
$ws.Range("F7").Value = 208
$ws.Range("F8").Value = 1077
$ws.Range("F9").Value = 277
$ws.Range("F11").Value = 352
$ws.Range("F12").Value = 627
$ws.Range("F14").Value = 495
$ws.Range("F15").Value = 136
$ws.Range("F17").Value = 157
$ws.Range("F18").Value = 830
$ws.Range("F20").Value = 522
$ws.Range("F22").Value = 10
$ws.Range("F23").Value = 308
$ws.Range("F24").Value = 209
$ws.Range("F27").Value = 584
$ws.Range("F28").Value = 964
$ws.Range("F29").Value = 6
$ws.Range("F31").Value = 234
$ws.Range("F32").Value = 1027
$ws.Range("F35").Value = 272

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1045
$ws.Range("F5").Value = 1045
$ws.Range("F8").Value = 226
$ws.Range("F14").Value = 587
$ws.Range("F15").Value = 92
$ws.Range("F17").Value = 970
$ws.Range("F26").Value = 3719
$ws.Range("F31").Value = 24
$ws.Range("F33").Value = 111

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2402
$ws.Range("F9").Value = 1244
$ws.Range("F10").Value = 330
$ws.Range("F11").Value = 90

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2402
$ws.Range("F8").Value = 1244
$ws.Range("F9").Value = 330
$ws.Range("F10").Value = 90
$ws.Range("F11").Value = 134
$ws.Range("F12").Value = 368
$ws.Range("F13").Value = 765
$ws.Range("F14").Value = 208
$ws.Range("F16").Value = 1077
$ws.Range("F17").Value = 277
$ws.Range("F18").Value = 352
$ws.Range("F19").Value = 627
$ws.Range("F20").Value = 1045
$ws.Range("F21").Value = 495
$ws.Range("F23").Value = 157
$ws.Range("F24").Value = 830
$ws.Range("F26").Value = 522
$ws.Range("F28").Value = 308
$ws.Range("F30").Value = 209
$ws.Range("F33").Value = 584
$ws.Range("F34").Value = 964
$ws.Range("F35").Value = 587
$ws.Range("F36").Value = 587
$ws.Range("F37").Value = 92
$ws.Range("F39").Value = 234
$ws.Range("F47").Value = 1027
$ws.Range("F51").Value = 272
